# BurnDownChart & LogBook edits
# - Update logged-hours in the LogBook table (rows 14-17)
# - Swap/rename two task descriptions (D17 / D18)
# - Update the saved selection on Sheet1 (was H21 -> now H19)
#
# Row 21/22 (the burn-down summary rows) and the embedded chart both read
# from these cells via SUM()/formula references, so they recompute
# automatically once the underlying log values below are changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- LogBook: hours logged per day for each task (rows 14-17) ---
$ws.Range("E14:J14").Value = 6
$ws.Range("E15:J15").Value = 4
$ws.Range("E16:J16").Value = 1
$ws.Range("E17:J17").Value = 1

# --- LogBook: task names for rows 17 & 18 are swapped, and the old
#     "Re-design tampilan" text is replaced with "Kegunaan Benchmark" ---
$ws.Range("D17").Value = "Kegunaan Benchmark"
$ws.Range("D18").Value = "Benchmark Apps"

# --- Update the sheet's saved selection / active cell ---
$ws.Range("H19").Select()

$wb.Save()
